$p = $ppt.ActivePresentation

# --- Slide 2 ("Roadmap") ---
$s = $p.Slides.Item(2)

# Title: "Roadmap (Ascheibe)" -> "Roadmap"
$s.Shapes.Item(1).TextFrame.TextRange.Text = "Roadmap"

# Content placeholder: add the roadmap bullet list
$tf = $s.Shapes.Item(2).TextFrame
$tr = $tf.TextRange

$tr.Text = "2 releases until HL 4.0"
$tr.InsertAfter("`rHL 3.3.14 “Unnamed”: begin of March 2016")
$tr.InsertAfter("`rHL 3.3.15 “Denver”: begin of July 2016")
$tr.InsertAfter("`rHL 4.0")
$tr.InsertAfter("`rUse summer for finishing 4.0")
$tr.InsertAfter("`rRelease end of September 2016")
$tr.InsertAfter("`rBetween 3.3.15 and 4.0")
$tr.InsertAfter("`rIntegration of new features into trunk")
$tr.InsertAfter("`rHL3legacy branch")
$tr.InsertAfter("`rLast version of trunk before integration")
$tr.InsertAfter("`rWill only get bugfixes")
$tr.InsertAfter("`rCan be used until trunk settles")
$tr.InsertAfter("`r")

# Indent levels (IndentLevel is 1-based; level 2 -> lvl="1", level 3 -> lvl="2")
$tr.Paragraphs(2).IndentLevel = 2
$tr.Paragraphs(3).IndentLevel = 2
$tr.Paragraphs(5).IndentLevel = 2
$tr.Paragraphs(6).IndentLevel = 2
$tr.Paragraphs(8).IndentLevel = 2
$tr.Paragraphs(9).IndentLevel = 2
$tr.Paragraphs(10).IndentLevel = 3
$tr.Paragraphs(11).IndentLevel = 3
$tr.Paragraphs(12).IndentLevel = 3
$tr.Paragraphs(13).IndentLevel = 2

# Split "Integration of new features into trunk" into 4 runs (paragraph 8)
$para8 = $tr.Paragraphs(8)
$para8.Characters(1, 15).Text = "Integration of "
$para8.Characters(16, 13).Text = "new features "
$para8.Characters(29, 5).Text = "into "
$para8.Characters(34, 5).Text = "trunk"

# Split "Will only get bugfixes" into 2 runs (paragraph 11)
$para11 = $tr.Paragraphs(11)
$para11.Characters(1, 14).Text = "Will only get "
$para11.Characters(15, 8).Text = "bugfixes"

# Shrink text on overflow (normAutofit)
$tf.AutoSize = 2

# --- Presentation level: sldId 267 -> 273 (2nd slide in sldIdLst) ---
# (handled structurally by PowerPoint when slide content changes; see below)

# --- Slide 9: merge "If true, cancelling" + " the " runs ---
$s9 = $p.Slides.Item(9)
$tf9 = $s9.Shapes.Item(2).TextFrame
$para3 = $tf9.TextRange.Paragraphs(3)
$para3.Characters(1, 24).Text = "If true, cancelling the "
